$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.050.48"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "2.299.01"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'300.21"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").Value = "'98.34"
$ws.Range("E6").Value = "  -1.54%  "
$ws.Range("D7").Value = "'0.520"
$ws.Range("E7").Value = "  +1.46%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.515"
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("D10").Value = "'36.22"
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").Value = "'0.0789"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D13").Value = "'17.71"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").Value = "'6.86"
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("D15").Value = "2.656.86"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "2.290.72"
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("D18").Value = "42.926.83"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "'12.81"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Value = "0.0₃0913"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").Value = "'68.98"
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("D23").Value = "'237.45"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").Value = "'2.13"
$ws.Range("E24").Value = "  -2.73%  "
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("D27").Value = "'24.91"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("D28").Value = "'165.33"
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").Value = "'9.11"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("D31").Value = "'33.04"
$ws.Range("E31").Value = "  -4.33%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").Value = "'5.08"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").Value = "'4.76"
$ws.Range("E34").Value = "  +2.04%  "
$ws.Range("D35").Value = "'17.91"
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").Value = "'0.0697"
$ws.Range("E37").Value = "  +0.76%  "
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("E39").Value = "  -0.80%  "
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").Value = "2.014.92"
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("D44").Value = "'2.24"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").Value = "'10.33"
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("D46").Value = "'17.48"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("D47").Value = "'2.83"
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("D48").Value = "'54.09"
$ws.Range("E48").Value = "  -2.59%  "
$ws.Range("D49").Value = "2.524.18"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("E50").Value = "  -1.56%  "
$ws.Range("D51").Value = "'73.25"
$ws.Range("E51").Value = "  +3.42%  "
